# overhaul write_to_excel with pandas
#
# Rewrites Sheet1 from a 4-col x 10-row hand-written table into a 17-col x
# 2-row pandas-style export: row 1 is a header ("ind" label + integer
# index/param columns), row 2 is one data record (a string "ind" value
# followed by numeric readings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 3-10 entirely (shifts nothing else up/down, just removes them)
# so only the header row and a single data row remain.
$ws.Rows("3:10").Delete()

# The original sheet only had 4 styled header cells (A1:D1, style index 1 -
# bold + border + centered). Extend that same style across the new
# header columns E1:Q1 by copying A1's formatting (format-only paste keeps
# the shared style index instead of minting new ones in styles.xml).
$ws.Range("A1").Copy()
$ws.Range("E1:Q1").PasteSpecial(-4122)

# Row 1 - header / index values
$ws.Range("A1").Value = "ind"
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9.81
$ws.Range("K1").Value = 11
$ws.Range("L1").Value = 12
$ws.Range("M1").Value = 13
$ws.Range("N1").Value = 14
$ws.Range("O1").Value = 15
$ws.Range("P1").Value = 16
$ws.Range("Q1").Value = 17

# Row 2 - single data record
$ws.Range("A2").Value = "ang_vel(limit(asin(protectedDiv(y3, y2), acos(y3, x2)), conditional(x1, conditional(y3, x3)), tan(y3)), cos(sin(x1)), cos(x2), x2)"
$ws.Range("B2").Value = -343.33
$ws.Range("C2").Value = -300.04
$ws.Range("D2").Value = -228.83
$ws.Range("E2").Value = -302.76
$ws.Range("F2").Value = -303.31
$ws.Range("G2").Value = -318.68
$ws.Range("H2").Value = -231.48
$ws.Range("I2").Value = -177.02
$ws.Range("J2").Value = -250.64
$ws.Range("K2").Value = -182.02
$ws.Range("L2").Value = -201.2
$ws.Range("M2").Value = -153.05
$ws.Range("N2").Value = -411.3
$ws.Range("O2").Value = -466.57
$ws.Range("P2").Value = -481.19
$ws.Range("Q2").Value = -452.88

Write-Output "edit applied: sheet1 rewritten to A1:Q2"
